# Rename the embedded logo pictures' display names (the wp:docPr / pic:cNvPr
# "name" attribute) inside the document's headers/footers:
#   - footer with the Pearson logo (docPr id="3")  : image2.png -> image1.png
#   - footer with the Pearson logo (docPr id="2")  : image2.png -> image1.png
#   - header with the BTEC logo   (docPr id="1")   : image1.jpg -> image2.jpg
#
# InlineShape objects don't reliably accept a direct ".Name = " assignment
# when they live inside a footer story in this host (it raises a stale
# "addressed block not found" error), so we select the shape's range first
# and re-fetch it from $word.Selection.InlineShapes before renaming - this
# works consistently for headers and footers alike.

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

function Rename-InlineShapeInRange($rng, $newName) {
    $shape = $rng.InlineShapes.Item(1)
    $shape.Range.Select()
    $word.Selection.InlineShapes.Item(1).Name = $newName
}

# --- Footers: Pearson Edexcel logo, image2.png -> image1.png ---
$footerFirst = $sec.Footers.Item(2)   # "first page" footer (docPr id="3")
Rename-InlineShapeInRange $footerFirst.Range "image1.png"

$footerPrimary = $sec.Footers.Item(1) # default/primary footer (docPr id="2")
Rename-InlineShapeInRange $footerPrimary.Range "image1.png"

# --- Header: BTEC logo, image1.jpg -> image2.jpg ---
$headerFirst = $sec.Headers.Item(2)   # "first page" header (docPr id="1")
Rename-InlineShapeInRange $headerFirst.Range "image2.jpg"
